$d = $word.ActiveDocument

# The body paragraph under the "Declaração do Problema" title reads (runs
# in order, "(b)" marks bold runs):
#
#   "O " "problema"(b) " referente a falta de controle sobre os data
#   centers " "afeta"(b) " os funcionários que trabalham nesses locais e
#   os sistemas da Telefônica " "devido"(b) " ao alto gasto de energia e
#   às quedas de energia, como desligamento dos servidores, instabilidade
#   e perda de dados."
#
# The edit swaps the order of the "afeta ..." and "devido ..." clauses
# and rewrites the trailing clause, producing:
#
#   "O " "problema"(b) " referente a falta de controle sobre os data
#   centers " "devido"(b) " ao alto gasto de energia e às quedas de
#   energia, como desligamento dos servidores, instabilidade e perda de
#   dados, isso " "afeta"(b) " os funcionários que trabalham nesses
#   locais e o dono da empresa."
#
# The four runs that change keep their original bold/plain formatting,
# so each one can be overwritten in place via Range.Text, which preserves
# the formatting of the run(s) it overwrites. Locate each run's Range via
# Find (so the script isn't dependent on hard-coded character offsets),
# then apply the replacements from the end of the paragraph towards the
# start so that not-yet-processed ranges keep valid Start/End positions.

function Get-MatchRange($searchText) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r
}

$rAfetaOld  = Get-MatchRange("afeta")
$rMidOld    = Get-MatchRange(" os funcionários que trabalham nesses locais e os sistemas da Telefônica ")
$rDevidoOld = Get-MatchRange("devido")
$rTailOld   = Get-MatchRange(" ao alto gasto de energia e às quedas de energia, como desligamento dos servidores, instabilidade e perda de dados.")

# Apply right-to-left (highest Start first) so earlier ranges stay valid.
$d.Range($rTailOld.Start, $rTailOld.End).Text = " os funcionários que trabalham nesses locais e o dono da empresa."
$d.Range($rDevidoOld.Start, $rDevidoOld.End).Text = "afeta"
$d.Range($rMidOld.Start, $rMidOld.End).Text = " ao alto gasto de energia e às quedas de energia, como desligamento dos servidores, instabilidade e perda de dados, isso "
$d.Range($rAfetaOld.Start, $rAfetaOld.End).Text = "devido"
